$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Add two new floating VML drawings (a flipped connector + a "ContextBox"
#    group) right before the run that hosts shape _x0000_s1044.
# ---------------------------------------------------------------------------

$newRunsXml = '<w:r><w:rPr><w:noProof/></w:rPr><w:pict><v:shape id="_x0000_s1111" type="#_x0000_t32" style="position:absolute;margin-left:454.1pt;margin-top:9.05pt;width:39.5pt;height:75.5pt;flip:x;z-index:251646967" o:connectortype="straight"/></w:pict></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:pict><v:group id="_x0000_s1112" style="position:absolute;margin-left:485.65pt;margin-top:-.25pt;width:81.35pt;height:45pt;z-index:251694080" coordorigin="2311,3945" coordsize="2385,900"><v:group id="_x0000_s1113" style="position:absolute;left:2311;top:3945;width:2385;height:900" coordorigin="270,7793" coordsize="2385,900"><v:rect id="_x0000_s1114" style="position:absolute;left:270;top:7793;width:2385;height:900"><v:textbox style="mso-next-textbox:#_x0000_s1114"><w:txbxContent><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>ContextBox</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect><v:shape id="_x0000_s1115" type="#_x0000_t32" style="position:absolute;left:270;top:8175;width:2385;height:0" o:connectortype="straight"/></v:group><v:shape id="_x0000_s1116" type="#_x0000_t32" style="position:absolute;left:2326;top:4590;width:2370;height:0" o:connectortype="straight"/></v:group></w:pict></w:r>'

$targetParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.WordOpenXML -like '*id="_x0000_s1044"*') {
        $targetParaIndex = $i
        break
    }
}

$insertAt = $d.Paragraphs.Item($targetParaIndex).Range
$insertAt.Collapse(1)
$wrapXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office"><w:body><w:p>' + $newRunsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertAt.InsertXML($wrapXml)

# ---------------------------------------------------------------------------
# 2) Remove the <w:proofErr .../> spell-check wrappers that used to surround
#    the ScoreBox / GameState / OptionsScreen textbox labels. Those nodes
#    live deep inside nested v:textbox/w:txbxContent content that isn't
#    individually addressable, so instead we round-trip the whole paragraph
#    that owns them through WordOpenXML (whose serializer drops proofErr
#    markers), delete the original paragraph, and splice the cleaned copy
#    back into the same slot.
# ---------------------------------------------------------------------------

$cleanParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.WordOpenXML -like '*<w:t>ScoreBox</w:t>*') {
        $cleanParaIndex = $i
        break
    }
}

$paraRange = $d.Paragraphs.Item($cleanParaIndex).Range
$capturedPackageXml = $paraRange.WordOpenXML

# WordOpenXML wraps a single paragraph's content as a tiny standalone
# document (renumbered paraId/rsid, plus a synthetic trailing empty
# paragraph + sectPr). Pull just the real paragraph's markup back out -
# its proofErr children are dropped by this serializer, which is exactly
# the cleanup we need.
$m = [regex]::Match($capturedPackageXml, '(?s)<w:body>(.*)<w:p[^>]*w:rsidR="00000000"[^>]*/><w:sectPr')
$cleanParaXml = $m.Groups[1].Value

$deleteStart = $paraRange.Start
$deleteEnd = $paraRange.End
$delRange = $d.Range($deleteStart, $deleteEnd)
$delRange.Delete()

$reinsertWrapXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' + $cleanParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$reinsertAt = $d.Range($deleteStart, $deleteStart)
$reinsertAt.InsertXML($reinsertWrapXml)
